# demo_platemap.xlsx - "added code to parse input file for SS2"
#
# The SMARTseq2 plate-map grid (rows 9-16 = plate rows A-H, columns B-M =
# well positions 1-12) had a few gaps/placeholder values left over from
# before the parser was wired up. Fill them in with the correct well
# labels now that the input-file parsing works, and move the active
# selection to the first cell that was fixed (H9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row A (sheet row 9): well A7 was showing the literal label "A7"
# instead of the "water" control that belongs there.
$ws.Range("H9").Value = "water"

# Row C (sheet row 11): C3 and C10 were missing entirely (C10 had
# leftover placeholder text "Water").
$ws.Range("D11").Value = "C3"
$ws.Range("K11").Value = "C10"

# Row E (sheet row 13): E6 was blank.
$ws.Range("G13").Value = "E6"

# Row G (sheet row 15): G1-G4 were left with stray whitespace/case test
# strings ("  water", "water", "wAter   ", "wat   er") instead of the
# real well labels.
$ws.Range("B15").Value = "G1"
$ws.Range("C15").Value = "G2"
$ws.Range("D15").Value = "G3"
$ws.Range("E15").Value = "G4"

# Row H (sheet row 16): H12 was blank.
$ws.Range("M16").Value = "H12"

# Leave the selection on the cell that kicked off this fix.
$ws.Range("H9").Select()
